# Network_Validation.xlsx — "Add files via upload" edit
#
# Content change: on the "Network_Diagnostic" sheet, the sample/placeholder
# value under the "Location" column (cell G2) changes from "eastus" to
# "location". Setting it directly also causes the shared-string table to
# drop the now-unused "eastus" entry and append "location" at the end,
# which is what shifts every other shared-string index in the workbook —
# matching the rest of the diff automatically.

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is active/selected right now so we can restore
# it — selecting a range on another sheet activates that sheet as a
# side-effect, and we don't want to change which tab is selected.
$originalActiveSheet = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("Network_Diagnostic")

# The actual content fix.
$ws.Range("G2").Value = "location"

# The sheet's view also moved: selection is now F8 (scrolled over to show
# the newly relevant columns).
$ws.Activate()
[void]$ws.Range("F8").Select()

# Restore the originally active/selected sheet (Check_Diskspace).
[void]$wb.Worksheets.Item($originalActiveSheet).Activate()
